$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "68.848.62"
$ws.Range("E2").Value = "  -0.30%  "

# Row 3
$ws.Range("D3").Value = "3.922.20"
$ws.Range("E3").Value = "  +3.00%  "

# Row 4
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.12%  "

# Row 5
$ws.Range("D5").Value = "'603.83"
$ws.Range("E5").Value = "  +0.49%  "

# Row 6
$ws.Range("D6").Value = "'167.74"
$ws.Range("E6").Value = "  +2.36%  "

# Row 7
$ws.Range("D7").Value = "3.916.85"
$ws.Range("E7").Value = "  +2.91%  "

# Row 8
$ws.Range("E8").Value = "  +0.03%  "

# Row 9
$ws.Range("E9").Value = "  -0.31%  "

# Row 10
$ws.Range("E10").Value = "  +0.34%  "

# Row 11
$ws.Range("E11").Value = "  +2.68%  "

# Row 12
$ws.Range("E12").Value = "  +0.99%  "

# Row 13
$ws.Range("D13").Value = "'0.0000255"
$ws.Range("E13").Value = "  +4.04%  "

# Row 14
$ws.Range("D14").Value = "'37.60"
$ws.Range("E14").Value = "  +1.30%  "

# Row 15
$ws.Range("D15").Value = "4.583.33"
$ws.Range("E15").Value = "  +3.18%  "

# Row 16
$ws.Range("D16").Value = "3.917.48"
$ws.Range("E16").Value = "  +2.74%  "

# Row 17
$ws.Range("D17").Value = "68.963.33"

# Row 18
$ws.Range("E18").Value = "  +0.36%  "

# Row 19
$ws.Range("D19").Value = "'17.46"
$ws.Range("E19").Value = "  +1.15%  "

# Row 20
$ws.Range("E20").Value = "  -2.00%  "

# Row 21
$ws.Range("E21").Value = "  -3.40%  "

# Row 22
$ws.Range("D22").Value = "'493.72"
$ws.Range("E22").Value = "  +1.58%  "

# Row 23
$ws.Range("D23").Value = "'0.730"
$ws.Range("E23").Value = "  +1.41%  "

# Row 24
$ws.Range("D24").Value = "'0.0000167"
$ws.Range("E24").Value = "  +4.11%  "

# Row 25
$ws.Range("D25").Value = "'84.89"
$ws.Range("E25").Value = "  +0.37%  "

# Row 26
$ws.Range("E26").Value = "  -0.02%  "

# Row 27
$ws.Range("D27").Value = "'12.13"
$ws.Range("E27").Value = "  -0.61%  "

# Row 28
$ws.Range("D28").Value = "'10.18"
$ws.Range("E28").Value = "  +1.31%  "

# Row 29
$ws.Range("E29").Value = "  +0.14%  "

# Row 30
$ws.Range("E30").Value = "  -0.11%  "

# Row 31
$ws.Range("D31").Value = "4.074.64"
$ws.Range("E31").Value = "  +2.81%  "

# Row 32
$ws.Range("D32").Value = "'2.38"
$ws.Range("E32").Value = "  -0.39%  "

# Row 33
$ws.Range("D33").Value = "'7.76"
$ws.Range("E33").Value = "  -3.15%  "

# Row 34
$ws.Range("D34").Value = "'32.00"
$ws.Range("E34").Value = "  +0.51%  "

# Row 35
$ws.Range("D35").Value = "3.882.83"
$ws.Range("E35").Value = "  +3.48%  "

# Row 36
$ws.Range("E36").Value = "  +0.32%  "

# Row 37
$ws.Range("E37").Value = "  +1.61%  "

# Row 38
$ws.Range("E38").Value = "  +0.68%  "

# Row 39
$ws.Range("E39").Value = "  +1.69%  "

# Row 40
$ws.Range("E40").Value = "  +7.35%  "

# Row 41
$ws.Range("D41").Value = "'1.00"
$ws.Range("E41").Value = "  +0.15%  "

# Row 42
$ws.Range("E42").Value = "  +0.15%  "

# Row 43
$ws.Range("D43").Value = "'437.07"
$ws.Range("E43").Value = "  -0.07%  "

# Row 44
$ws.Range("D44").Value = "'1.99"
$ws.Range("E44").Value = "  +0.58%  "

# Row 45
$ws.Range("D45").Value = "'48.03"
$ws.Range("E45").Value = "  -1.23%  "

# Row 48
$ws.Range("D48").Value = "'143.05"
$ws.Range("E48").Value = "  +0.64%  "

# Row 49
$ws.Range("B49").Value = "Maker"
$ws.Range("C49").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D49").Value = "2.823.68"
$ws.Range("E49").Value = "  -0.14%  "

# Row 50
$ws.Range("B50").Value = "FLOKI"
$ws.Range("C50").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D50").Value = "'0.000269"
$ws.Range("E50").Value = "  +18.12%  "

# Row 51
$ws.Range("B51").Value = "VeChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D51").Value = "'0.0358"
$ws.Range("E51").Value = "  +1.82%  "
